$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet was renamed as part of unifying the DataNode/DataTable/Entity
# concepts: "Property1" -> "DataNode".
$ws.Name = "DataNode"

# The author's selection moved from K17 to B41 before saving.
$ws.Range("B41").Select()
